# Update column G ("K") values for rows 2-23 on the active sheet,
# reflecting the regenerated save_data (K instead of Strike#, std/mean,
# s_vals recomputed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 7
    3  = 5
    4  = 5
    5  = 4
    6  = 4
    7  = 3
    8  = 5
    9  = 5
    10 = 1
    11 = 4
    12 = 1
    13 = 1
    14 = 2
    15 = 2
    16 = 2
    17 = 1
    18 = 5
    19 = 0
    20 = 2
    21 = 3
    22 = 2
    23 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
